# Update the "Förändrad" (Changed) date column from 2023-09-10 (45179)
# to 2023-09-11 (45180) for all data rows (C2:C150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C150")
$range.Value = 45180
